# Updated Results with corrected code
# Applies corrected figures to Sheet1 of 2040_EL.xlsx:
#  - D3 (Hydrogen / Non-metallic minerals) is cleared
#  - C4 (Methanol / Chemicals) corrected to 0
#  - C5 (Ammonia / Chemicals) corrected
#  - Row 7 label renamed "Other" -> "Biogas", D7 value corrected
#  - A new row 8 "Other" is appended with a corrected D8 value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray Hydrogen / Non-metallic minerals figure
$ws.Range("D3").Value = ""

# Corrected Methanol / Chemicals figure
$ws.Range("C4").Value = 0

# Corrected Ammonia / Chemicals figure
$ws.Range("C5").Value = 585.7293526225917

# Row 7 is now "Biogas" (was "Other") with a corrected value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 248.1905827916965

# New row 8: re-introduce "Other" below Biogas, matching the header row's
# label formatting (bold, centered, bordered) used by the other labels in
# column A.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 3842.879822249616
